# Insert a new data row above row 135 (pushing existing rows 135..205 down to 136..206)
# and populate the new row 135 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(135).Insert()

$ws.Cells.Item(135, 1).Value2 = 10
$ws.Cells.Item(135, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(135, 3).Value2 = "La Araucanía"
$ws.Cells.Item(135, 4).Value2 = 44488
$ws.Cells.Item(135, 5).Value2 = 9
$ws.Cells.Item(135, 6).Value2 = 100112009
$ws.Cells.Item(135, 7).Value2 = "Acelga"
$ws.Cells.Item(135, 8).Value2 = "Sin especificar"
$ws.Cells.Item(135, 9).Value2 = "Primera"
$ws.Cells.Item(135, 10).Value2 = 120
$ws.Cells.Item(135, 11).Value2 = 7000
$ws.Cells.Item(135, 12).Value2 = 8000
$ws.Cells.Item(135, 13).Value2 = 7542
$ws.Cells.Item(135, 14).Value2 = "$/docena de atados (12 kilos)"
$ws.Cells.Item(135, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(135, 16).Value2 = 628
$ws.Cells.Item(135, 17).Value2 = 12
$ws.Cells.Item(135, 18).Value2 = "Hortaliza"
